# Applies the "Updated pom with rest-api-version" change:
#  - Populates the (previously empty) "Sheet1" worksheet with a new
#    config-upload test row, cloning the layout/styles already used by the
#    "FileUpload" worksheet.
#  - Adds 3 new shared strings used by that new row (a URL, a JSON payload,
#    and a curl-style upload description), and wires a hyperlink on B3 to
#    the URL.
#  - Nudges the saved cursor/selection on the "FileUpload" and
#    "SuiteVariable" worksheets to F3.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("FileUpload")
$ws3 = $wb.Worksheets.Item("SuiteVariable")

# ---------------------------------------------------------------------
# 1) FileUpload: move the saved viewport/selection.
# ---------------------------------------------------------------------
$ws2.Range("F3").Select()

# ---------------------------------------------------------------------
# 2) SuiteVariable: move the saved selection.
# ---------------------------------------------------------------------
$ws3.Range("F3").Select()

# ---------------------------------------------------------------------
# 3) Sheet1: build the new table, mirroring FileUpload's layout.
# ---------------------------------------------------------------------

# -- column widths (approximate "characters" width; engine quantizes to px)
$ws1.Columns.Item(1).ColumnWidth = 8.43
$ws1.Columns.Item(2).ColumnWidth = 37.71
$ws1.Columns.Item(4).ColumnWidth = 12
$ws1.Columns.Item(5).ColumnWidth = 8
$ws1.Columns.Item(6).ColumnWidth = 89.71
$ws1.Columns.Item(7).ColumnWidth = 12.43
$ws1.Columns.Item(8).ColumnWidth = 9.86
$ws1.Columns.Item(12).ColumnWidth = 8.71
$ws1.Columns.Item(14).ColumnWidth = 11.14

# -- row 1 (header) : copy formatting + values from FileUpload row 1
$row1Cols = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14)
foreach ($col in $row1Cols) {
    $src = $ws2.Cells.Item(1, $col)
    $dst = $ws1.Cells.Item(1, $col)
    $src.Copy()
    $dst.PasteSpecial(-4122) | Out-Null
}
$ws1.Cells.Item(1,1).Value2 = $ws2.Cells.Item(1,1).Value2
$ws1.Cells.Item(1,2).Value2 = $ws2.Cells.Item(1,2).Value2
$ws1.Cells.Item(1,3).Value2 = $ws2.Cells.Item(1,3).Value2
$ws1.Cells.Item(1,4).Value2 = $ws2.Cells.Item(1,4).Value2
$ws1.Cells.Item(1,7).Value2 = $ws2.Cells.Item(1,7).Value2
$ws1.Cells.Item(1,13).Value2 = $ws2.Cells.Item(1,13).Value2
$ws1.Cells.Item(1,14).Value2 = $ws2.Cells.Item(1,14).Value2

# -- row 2 (sub-header) : only columns D..L are populated
$row2Cols = @(4,5,6,7,8,9,10,11,12)
foreach ($col in $row2Cols) {
    $src = $ws2.Cells.Item(2, $col)
    $dst = $ws1.Cells.Item(2, $col)
    $src.Copy()
    $dst.PasteSpecial(-4122) | Out-Null
    $dst.Value2 = $src.Value2
}
$ws1.Rows.Item(2).RowHeight = 60

# -- row 3 (data row)
# A3, C3 reuse the same shared strings as FileUpload's A3/C3
foreach ($col in @(1,3)) {
    $src = $ws2.Cells.Item(3, $col)
    $dst = $ws1.Cells.Item(3, $col)
    $src.Copy()
    $dst.PasteSpecial(-4122) | Out-Null
    $dst.Value2 = $src.Value2
}

# B3 / E3 / F3 carry new content - copy formatting only, then set new text
foreach ($col in @(2,5,6)) {
    $src = $ws2.Cells.Item(3, $col)
    $dst = $ws1.Cells.Item(3, $col)
    $src.Copy()
    $dst.PasteSpecial(-4122) | Out-Null
}

$ws1.Cells.Item(3,2).Value2 = "http://10.0.3.41:8080/rulify/v1/config_upload/OSLOS-1116-xlsx2"
$ws1.Cells.Item(3,5).Value2 = "{`n   ""APPLICATION_ID"":""1"",`n   ""TENANT_ID"":""1"",`n   ""UID"":""abcd""`n}"
$ws1.Cells.Item(3,6).Value2 = "F file=@C:/Users/dwiveddi/Desktop/ActionDriverAPI/OSLOS-1116-xlsx2.xlsx`nF description=cvd"

$ws1.Cells.Item(3,7).Value2 = 200

$ws1.Rows.Item(3).RowHeight = 195

# -- merged header cells
$ws1.Range("D1:F1").Merge() | Out-Null
$ws1.Range("G1:I1").Merge() | Out-Null
$ws1.Range("J1:L1").Merge() | Out-Null

# -- hyperlink on B3, pointing at the same URL as its display text
$ws1.Hyperlinks.Add($ws1.Range("B3"), "http://10.0.3.41:8080/rulify/v1/config_upload/OSLOS-1116-xlsx2") | Out-Null

# -- selection / active sheet : Sheet1 stays the selected tab
$ws1.Activate()
$ws1.Range("F3").Select()
